$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Activités")
$ws2 = $wb.Worksheets.Item("Données")

# ---------------------------------------------------------------
# "Activités" sheet - fill in the journal entries for rows 18-20
# ---------------------------------------------------------------

# Row 18 - finish the entry that only had a start time
$ws1.Range("C18").Value = 0.51041666666666663
$ws1.Range("E18").Value = "Mise en place de l'infrastructure"
$ws1.Range("F18").Value = "Analyse"
$ws1.Range("G18").Value = "Changement de Github"

# Row 19 - new entry
$ws1.Range("A19").Value = 44259
$ws1.Range("B19").Value = 0.57916666666666672
$ws1.Range("C19").Value = 0.62847222222222221
$ws1.Range("E19").Value = "Création"
$ws1.Range("F19").Value = "Analyse"
$ws1.Range("G19").Value = "Création d'un projet de test"

# Row 20 - new entry
$ws1.Range("A20").Value = 44260
$ws1.Range("B20").Value = 0.44444444444444442
$ws1.Range("C20").Value = 0.62430555555555556
$ws1.Range("E20").Value = "Création"
$ws1.Range("F20").Value = "Analyse"
$ws1.Range("G20").Value = "Création du projet principal et création de l'interface"

# ---------------------------------------------------------------
# "Données" sheet - register the new activity in the lookup list
# ---------------------------------------------------------------
$ws2.Range("A9").Value = "Création"
$ws2.Range("B9").Value = "Création"

# ---------------------------------------------------------------
# Selections - match the cursor position left by the author
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A9").Select()

$ws1.Activate()
$ws1.Range("G16").Select()
